$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 3778808.76
$ws.Range("C7").Value = -14.95065677575798
$ws.Range("D7").Value = 3373
$ws.Range("E7").Value = 3373
$ws.Range("F7").Value = 1120.310927957308
$ws.Range("G7").Value = 19.41704402905728
